$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 38: ammo_12x76_zhekan (Slug) - increase DMG (H38) from 2.5 to 2.7
$ws.Range("H38").Value = 2.7

# Row 39: ammo_12x70_buck (Buckshot) - increase DMG formula (H39) from 9*0.4 to 9*0.42
$ws.Range("H39").Formula = "=9*0.42"

# Update selection to match authored state
$ws.Range("J27").Select()
